$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "246.02"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.085"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05611"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.477"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.020"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8115"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8478"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1342"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.03229"
$ws.Range("E11").Value = "10LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.02774"
$ws.Range("E12").Value = "11BitrueCoinBTR"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09406"
$ws.Range("E13").Value = "12BitMartTokenBMX"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001512"
$ws.Range("E14").Value = "13BitForexTokenBF"
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006019"
$ws.Range("E15").Value = "14OneONEWorstin24h"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006137"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.557"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06995"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1320"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.738"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04688"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001249"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004619"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.00009597"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001389"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03656"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006108"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002499"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008621"
$ws.Range("E44").Value = "43LocalTradersLCTBestin24h"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005294"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002099"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0001999"
